$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I2:I46 and J2:J46 with the new data values
$iValues = @(
    8,7,8,7,7,6,5,7,1,9,6,7,5,6,7,6,7,9,9,7,6,7,7,8,7,6,1,8,5,1,5,1,9,6,7,12,4,9,9,7,7,8,6,5,5
)
$jValues = @(
    8,7,8,7,7,7,5,7,2,9,6,8,5,7,7,6,7,9,9,7,7,7,7,8,7,8,2,8,6,2,6,2,9,6,7,13,5,10,9,8,7,8,6,5,5
)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}

